# Update the "想去人数" (F column) counts across the sheets, per the
# published diff. Each entry is (Row, NewValue).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ 4 = 241; 5 = 1976; 6 = 219; 7 = 644; 8 = 15; 9 = 134; 11 = 631; 12 = 22; 13 = 74; 14 = 873 }
    "演出"     = @{ 6 = 8; 11 = 25; 17 = 32 }
    "本地生活" = @{ 5 = 179 }
    "全部类型" = @{ 6 = 179; 12 = 241; 13 = 8; 16 = 1976; 18 = 219; 20 = 25; 21 = 644; 22 = 15; 23 = 134; 26 = 631; 27 = 22; 28 = 74; 30 = 873; 37 = 32 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
